# Update "想去人数" (want-to-go count) figures in column F across sheets
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) to reflect
# the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1376
$ws1.Range("F3").Value  = 1765
$ws1.Range("F4").Value  = 906
$ws1.Range("F7").Value  = 682
$ws1.Range("F11").Value = 2489
$ws1.Range("F12").Value = 1609
$ws1.Range("F13").Value = 1524
$ws1.Range("F16").Value = 624
$ws1.Range("F19").Value = 316
$ws1.Range("F20").Value = 1092
$ws1.Range("F24").Value = 5167
$ws1.Range("F26").Value = 597
$ws1.Range("F27").Value = 89
$ws1.Range("F28").Value = 163
$ws1.Range("F29").Value = 139
$ws1.Range("F30").Value = 231
$ws1.Range("F32").Value = 34
$ws1.Range("F33").Value = 1046
$ws1.Range("F34").Value = 752
$ws1.Range("F36").Value = 58
$ws1.Range("F38").Value = 399
$ws1.Range("F39").Value = 1105
$ws1.Range("F40").Value = 134
$ws1.Range("F41").Value = 107
$ws1.Range("F42").Value = 178
$ws1.Range("F44").Value = 58

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 788
$ws2.Range("F6").Value = 9

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1376
$ws4.Range("F4").Value  = 788
$ws4.Range("F5").Value  = 1766
$ws4.Range("F6").Value  = 906
$ws4.Range("F11").Value = 682
$ws4.Range("F13").Value = 9
$ws4.Range("F17").Value = 2489
$ws4.Range("F18").Value = 1609
$ws4.Range("F19").Value = 1524
$ws4.Range("F22").Value = 624
$ws4.Range("F26").Value = 316
$ws4.Range("F27").Value = 1092
$ws4.Range("F29").Value = 5168
$ws4.Range("F31").Value = 597
$ws4.Range("F32").Value = 89
$ws4.Range("F33").Value = 163
$ws4.Range("F34").Value = 139
$ws4.Range("F35").Value = 231
$ws4.Range("F37").Value = 34
$ws4.Range("F38").Value = 1046
$ws4.Range("F39").Value = 752
$ws4.Range("F40").Value = 58
$ws4.Range("F41").Value = 399
$ws4.Range("F42").Value = 1105
$ws4.Range("F43").Value = 134
$ws4.Range("F44").Value = 178
$ws4.Range("F46").Value = 58
